# TFS6511 - eCL Pilot Survey Question - changed order of responses and added follow up on line 7
# TFS9511 - ecl survey pilot

$wb = $excel.ActiveWorkbook

$wsHistory = $wb.Worksheets.Item("Change History")
$wsSurvey = $wb.Worksheets.Item("eCL Survey")

# ---------------------------------------------------------------------------
# Sheet: eCL Survey (question 6 - "How prepared was your supervisor ...")
#   - Reverse the order of the 1-5 options (was Very Prepared -> Very
#     Unprepared, now Very Unprepared -> Very Prepared)
#   - Add a follow up prompt on that row (was "N/A", now "Please explain
#     below.")
# ---------------------------------------------------------------------------
$wsSurvey.Range("F7").Value = "1 - Very Unprepared`n2 - Unprepared`n3 - Neither Prepared or Unprepared`n4 - Prepared`n5 - Very Prepared"
$wsSurvey.Range("G7").Value = "Please explain below."

# ---------------------------------------------------------------------------
# Sheet: Change History - log this change as a new row
# ---------------------------------------------------------------------------
$wsHistory.Range("B9").Value = "1/29/2018"
$wsHistory.Range("C9").Value = "TFS6511 - eCL Pilot Survey Question - changed order of responses and added follow up on line 7"
$wsHistory.Range("D9").Value = 1.04
$wsHistory.Range("E9").Value = "Doug Stearns"
$wsHistory.Rows.Item(9).RowHeight = 39.6

# ---------------------------------------------------------------------------
# Restore the selection on each sheet, and leave "eCL Survey" as the active
# tab (as it was before), by selecting it last.
# ---------------------------------------------------------------------------
$wsHistory.Range("E9").Select()
$wsSurvey.Range("G7").Select()
